$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3584.2104
$ws.Range("I8").Value = 1911.1111
$ws.Range("K8").Value = 5733.3333
$ws.Range("M8").Value = -5594.3333

$ws.Range("H58").Value = 10359.8
$ws.Range("I58").Value = 600
$ws.Range("K58").Value = 1800
$ws.Range("M58").Value = -1650

$ws.Range("H64").Value = 4823.5293
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752

$ws.Range("H67").Value = 4823.5293
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142

$ws.Range("H81").Value = 76357.25
$ws.Range("J81").Value = 76357.25
$ws.Range("L81").Value = 76357.25
$ws.Range("N81").Value = -78353.25

$ws.Range("H84").Value = 76357.25
$ws.Range("J84").Value = 76357.25
$ws.Range("L84").Value = 229071.75
$ws.Range("N84").Value = -239055.75

$ws.Range("H107").Value = 2652.5
$ws.Range("I107").Value = 1791.375
$ws.Range("K107").Value = 1791.375
$ws.Range("M107").Value = 128.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 772.1667
$ws.Range("I4").Value = 334
$ws.Range("J4").Value = 1648.5
$ws.Range("K4").Value = 334
$ws.Range("L4").Value = 1648.5
$ws.Range("M4").Value = -218
$ws.Range("N4").Value = -1880.5

$ws.Range("H32").Value = 10874286
$ws.Range("I32").Value = 11368565
$ws.Range("K32").Value = 11368565
$ws.Range("M32").Value = -11368278

$ws.Range("H45").Value = 1999.8889
$ws.Range("I45").Value = 1215.4286
$ws.Range("K45").Value = 1215.4286
$ws.Range("M45").Value = -838.4286

$ws.Range("H110").Value = 8015.727
$ws.Range("I110").Value = 6034.125
$ws.Range("J110").Value = 13300
$ws.Range("K110").Value = 6034.125
$ws.Range("L110").Value = 13300
$ws.Range("M110").Value = -3989.125
$ws.Range("N110").Value = -17390

$ws.Range("H122").Value = 1940.2
$ws.Range("I122").Value = 1599.5
$ws.Range("K122").Value = 4798.5
$ws.Range("M122").Value = -2348.5

$ws.Range("H132").Value = 4661.778
$ws.Range("I132").Value = 2882.8333
$ws.Range("K132").Value = 8648.499899999999
$ws.Range("M132").Value = -6118.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1957.8334
$ws.Range("I105").Value = 1769.4
$ws.Range("K105").Value = 1769.4
$ws.Range("M105").Value = -22.40000000000009

$ws.Range("H107").Value = 1630.6666
$ws.Range("I107").Value = 1419.2307
$ws.Range("K107").Value = 1419.2307
$ws.Range("M107").Value = 500.7692999999999

$ws.Range("H134").Value = 102348
$ws.Range("I134").Value = 9476.75
$ws.Range("K134").Value = 28430.25
$ws.Range("M134").Value = -25895.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 575
$ws.Range("I16").Value = 575
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 575
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -288
$ws.Range("N16").ClearContents()

$ws.Range("H58").Value = 2451.9092
$ws.Range("I58").Value = 2421.75
$ws.Range("J58").Value = 2469.1428
$ws.Range("K58").Value = 2421.75
$ws.Range("L58").Value = 2469.1428
$ws.Range("M58").Value = -2218.75
$ws.Range("N58").Value = -2875.1428

$ws.Range("H113").Value = 575
$ws.Range("I113").Value = 575
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 575
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1595
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 3525.9473
$ws.Range("I122").Value = 1691.8462
$ws.Range("K122").Value = 5075.5386
$ws.Range("M122").Value = -2625.5386

$ws.Range("H134").Value = 838439
$ws.Range("I134").Value = 1668212.4
$ws.Range("K134").Value = 5004637.199999999
$ws.Range("M134").Value = -5002102.199999999

$ws.Range("H136").Value = 2451.9092
$ws.Range("I136").Value = 2421.75
$ws.Range("J136").Value = 2469.1428
$ws.Range("K136").Value = 7265.25
$ws.Range("L136").Value = 7407.428400000001
$ws.Range("M136").Value = -4715.25
$ws.Range("N136").Value = -12507.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 297.41666
$ws.Range("I2").Value = 414.625
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 2487.75
$ws.Range("L2").Value = 378
$ws.Range("M2").Value = -2374.75
$ws.Range("N2").Value = -604

$ws.Range("H23").Value = 859.4
$ws.Range("I23").Value = 824.25
$ws.Range("K23").Value = 2472.75
$ws.Range("M23").Value = -2237.75

$ws.Range("H34").Value = 2594
$ws.Range("J34").Value = 4999.5
$ws.Range("L34").Value = 14998.5
$ws.Range("N34").Value = -15166.5

$ws.Range("H116").Value = 782.25
$ws.Range("I116").Value = 543
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 1629
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = 1813
$ws.Range("N116").Value = -11384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 220.66667
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 333.2
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 333.2
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = -559.2

$ws.Range("H14").Value = 7500295
$ws.Range("I14").Value = 15000352
$ws.Range("K14").Value = 15000352
$ws.Range("M14").Value = -15000184

$ws.Range("H132").Value = 100001620
$ws.Range("I132").Value = 111112744
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 333338232
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -333335702
$ws.Range("N132").Value = -9557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1960.1428
$ws.Range("I22").Value = 2579.8
$ws.Range("K22").Value = 2579.8
$ws.Range("M22").Value = -2284.8

$ws.Range("H27").Value = 1960.1428
$ws.Range("I27").Value = 2579.8
$ws.Range("K27").Value = 2579.8
$ws.Range("M27").Value = -2472.8

$ws.Range("H40").Value = 3505.7778
$ws.Range("I40").Value = 1778.6666
$ws.Range("K40").Value = 1778.6666
$ws.Range("M40").Value = -1642.6666

$ws.Range("H46").Value = 1938.0555
$ws.Range("I46").Value = 1553.871
$ws.Range("K46").Value = 1553.871
$ws.Range("M46").Value = -1365.871

$ws.Range("H61").Value = 1375
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H113").Value = 1375
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 870.7143
$ws.Range("J81").Value = 599
$ws.Range("L81").Value = 1198
$ws.Range("N81").Value = -3320

$ws.Range("H84").Value = 870.7143
$ws.Range("J84").Value = 599
$ws.Range("L84").Value = 5990
$ws.Range("N84").Value = -16598

$ws.Range("H102").Value = 109995
$ws.Range("J102").Value = 109995
$ws.Range("L102").Value = 109995
$ws.Range("N102").Value = -116485

$ws.Range("H109").Value = 56998.75
$ws.Range("J109").Value = 56998.75
$ws.Range("L109").Value = 56998.75
$ws.Range("N109").Value = -59772.75

$ws.Range("H136").Value = 32999.875
$ws.Range("J136").Value = 43333.332
$ws.Range("L136").Value = 129999.996
$ws.Range("N136").Value = -135099.996
